$wb = $excel.ActiveWorkbook

# This script reapplies a scheduled market-data refresh (Universalis price pull)
# across the Golem_Profits workbook: per-leve currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) are overwritten with freshly fetched values. Rows whose
# profit no longer applies (e.g. NQ/HQ price now equal, or no craftable profit) have
# their M/N profit cells cleared entirely, matching how the exporter omits cells that
# hold no value.

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 19
$ws.Range("H19").Value = 777.8333
$ws.Range("I19").Value = 400.33334
$ws.Range("J19").Value = 903.6667
$ws.Range("K19").Value = 400.33334
$ws.Range("L19").Value = 903.6667
$ws.Range("M19").Value = -225.33334
$ws.Range("N19").Value = -1253.6667

# Row 40
$ws.Range("H40").Value = 4510.8887
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4510.8887
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4510.8887
$ws.Range("N40").Value = -4860.8887
$ws.Range("M40").ClearContents()

# Row 51
$ws.Range("H51").Value = 239949.75
$ws.Range("I51").Value = 205000
$ws.Range("K51").Value = 205000
$ws.Range("M51").Value = -204516

# Row 86
$ws.Range("H86").Value = 55666.668
$ws.Range("J86").Value = 55666.668
$ws.Range("L86").Value = 55666.668
$ws.Range("N86").Value = -57912.668

# Row 89
$ws.Range("H89").Value = 55666.668
$ws.Range("J89").Value = 55666.668
$ws.Range("L89").Value = 278333.34
$ws.Range("N89").Value = -289565.34

# Row 109
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774

# Row 138
$ws.Range("H138").Value = 2687.6191
$ws.Range("I138").Value = 1250
$ws.Range("K138").Value = 3750
$ws.Range("M138").Value = 1390

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 26
$ws.Range("H26").Value = 925
$ws.Range("I26").Value = 880
$ws.Range("K26").Value = 880
$ws.Range("M26").Value = -550

# Row 39
$ws.Range("H39").Value = 5603.2
$ws.Range("I39").Value = 4504
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 4504
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = -3984
$ws.Range("N39").Value = -11040

# Row 50
$ws.Range("H50").Value = 15662.333
$ws.Range("J50").Value = 43785
$ws.Range("L50").Value = 43785
$ws.Range("N50").Value = -45213

# Row 122
$ws.Range("H122").Value = 2828
$ws.Range("I122").Value = 2828
$ws.Range("K122").Value = 8484
$ws.Range("M122").Value = -6034

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 7
$ws.Range("H7").Value = 5043.375
$ws.Range("I7").Value = 28.6
$ws.Range("K7").Value = 28.6
$ws.Range("M7").Value = 84.40000000000001

# Row 134
$ws.Range("H134").Value = 383.33334
$ws.Range("I134").Value = 383.33334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 1150.00002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 1384.99998
$ws.Range("N134").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 31
$ws.Range("H31").Value = 1500
$ws.Range("I31").Value = 1500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1205

# Row 34
$ws.Range("H34").Value = 1500
$ws.Range("I34").Value = 1500
$ws.Range("K34").Value = 1500
$ws.Range("M34").Value = -1298

# Row 35
$ws.Range("H35").Value = 2537.75
$ws.Range("I35").Value = 2810.5715
$ws.Range("J35").Value = 628
$ws.Range("K35").Value = 2810.5715
$ws.Range("L35").Value = 628
$ws.Range("M35").Value = -2516.5715
$ws.Range("N35").Value = -1216

# Row 48
$ws.Range("H48").Value = 24665.666
$ws.Range("J48").Value = 24665.666
$ws.Range("L48").Value = 24665.666
$ws.Range("N48").Value = -25617.666

# Row 92
$ws.Range("H92").Value = 9000
$ws.Range("J92").Value = 9000
$ws.Range("L92").Value = 9000
$ws.Range("N92").Value = -13992

# Row 106
$ws.Range("H106").Value = 69021.125
$ws.Range("J106").Value = 69021.125
$ws.Range("L106").Value = 69021.125
$ws.Range("N106").Value = -71545.125

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 7
$ws.Range("H7").Value = 400.75
$ws.Range("I7").Value = 367
$ws.Range("K7").Value = 1101
$ws.Range("M7").Value = -989

# Row 40
$ws.Range("H40").Value = 456.57144
$ws.Range("I40").Value = 456.57144
$ws.Range("K40").Value = 1826.28576
$ws.Range("M40").Value = -1757.28576

# Row 63
$ws.Range("H63").Value = 4568.778

# Row 66
$ws.Range("H66").Value = 4568.778

# Row 103
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

# Row 131
$ws.Range("H131").Value = 3380.8572
$ws.Range("I131").Value = 3694.75
$ws.Range("J131").Value = 2962.3333
$ws.Range("K131").Value = 11084.25
$ws.Range("L131").Value = 8886.999899999999
$ws.Range("M131").Value = -6044.25
$ws.Range("N131").Value = -18966.9999

# Row 139
$ws.Range("H139").Value = 1028.7
$ws.Range("I139").Value = 1028.2222
$ws.Range("J139").Value = 1033
$ws.Range("K139").Value = 3084.6666
$ws.Range("L139").Value = 3099
$ws.Range("M139").Value = 2055.3334
$ws.Range("N139").Value = -13379

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# Row 43
$ws.Range("H43").Value = 14959.4
$ws.Range("I43").Value = 4966.6665
$ws.Range("J43").Value = 29948.5
$ws.Range("K43").Value = 4966.6665
$ws.Range("L43").Value = 29948.5
$ws.Range("M43").Value = -4815.6665
$ws.Range("N43").Value = -30250.5

# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

# Row 55
$ws.Range("H55").Value = 35666.668
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

# Row 93
$ws.Range("H93").Value = 49990
$ws.Range("J93").Value = 49990
$ws.Range("L93").Value = 49990
$ws.Range("N93").Value = -53734

# Row 95
$ws.Range("H95").Value = 26666.334
$ws.Range("J95").Value = 26666.334
$ws.Range("L95").Value = 26666.334
$ws.Range("N95").Value = -32158.334

# Row 101
$ws.Range("H101").Value = 9955
$ws.Range("J101").Value = 9955
$ws.Range("L101").Value = 9955
$ws.Range("N101").Value = -16445

# Row 113
$ws.Range("H113").Value = 1497
$ws.Range("I113").Value = 1497
$ws.Range("K113").Value = 1497
$ws.Range("M113").Value = 673

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 46
$ws.Range("H46").Value = 339332.34
$ws.Range("I46").Value = 668666.7
$ws.Range("J46").Value = 9998
$ws.Range("K46").Value = 668666.7
$ws.Range("L46").Value = 9998
$ws.Range("M46").Value = -668478.7
$ws.Range("N46").Value = -10374

# Row 61
$ws.Range("H61").Value = 6247
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 113
$ws.Range("H113").Value = 6247
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 55
$ws.Range("H55").Value = 25609.4
$ws.Range("I55").Value = 12524
$ws.Range("J55").Value = 34333
$ws.Range("K55").Value = 12524
$ws.Range("L55").Value = 34333
$ws.Range("M55").Value = -12247
$ws.Range("N55").Value = -34887

# Row 100
$ws.Range("H100").Value = 161.75
$ws.Range("I100").Value = 149
$ws.Range("J100").Value = 200
$ws.Range("K100").Value = 298
$ws.Range("L100").Value = 400
$ws.Range("M100").Value = 243
$ws.Range("N100").Value = -1482

# Row 124
$ws.Range("H124").Value = 32499.5
$ws.Range("I124").Value = 14999.5
$ws.Range("J124").Value = 49999.5
$ws.Range("K124").Value = 14999.5
$ws.Range("L124").Value = 49999.5
$ws.Range("M124").Value = -10089.5
$ws.Range("N124").Value = -59819.5
